$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.12950000000001
$ws.Range("B4").Value = 4.496400000000002
$ws.Range("C4").Value = -10.84139999999999

$ws.Range("B5").Value = 5.124399999999998

$ws.Range("A7").Value = -21.48820000000001

$ws.Range("B8").Value = 4.890300000000001

$ws.Range("C9").Value = -11.76230000000001

$ws.Range("A16").Value = -21.49730000000002
$ws.Range("B16").Value = 4.827500000000004

$ws.Range("C18").Value = -14.6768
